# "Add cantrals by cantons" — rebuild Sheet1 as a clean one-row-header table:
#   idx | idx2 | Name | Date Start | Date End | (m3/s) | (MW1) | (MW2) | (GWh) Winter | (GWh) Summer | (GWh) Year
# The old sheet had a 2-row header and data in rows 3-9 (with a stray,
# unnamed row 5). We snapshot the existing per-plant data (exact values,
# incl. float bit patterns) before touching anything, drop the orphan row,
# then rewrite the sheet with the new header + the same data shifted up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Snapshot the six named power-plant rows (old rows 3,4,6,7,8,9) ----
# Row 5 (F5=3 G5=5.3 H5=5.3, no name) is an orphan and is dropped.
$srcRows = @(3, 4, 6, 7, 8, 9)
$data = @()
foreach ($r in $srcRows) {
    $row = @{
        A = $ws.Cells.Item($r, 1).Value2
        B = $ws.Cells.Item($r, 2).Value2
        C = $ws.Cells.Item($r, 3).Value2
        D = $ws.Cells.Item($r, 4).Value2
        E = $ws.Cells.Item($r, 5).Value2
        F = $ws.Cells.Item($r, 6).Value2
        G = $ws.Cells.Item($r, 7).Value2
        H = $ws.Cells.Item($r, 8).Value2
        I = $ws.Cells.Item($r, 9).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
    }
    $data += $row
}

# ---- Clean slate ----
$ws.Range("A1:K9").Clear()

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Font.Name = "Arial"
    $c.Font.Size = 9
    $c.Value = $text
}

function Set-IntCell($row, $col, $n) {
    if ($n -eq $null) { return }
    $c = $ws.Cells.Item($row, $col)
    $c.Font.Name = "Arial"
    $c.Font.Size = 9
    $c.NumberFormat = "0"
    $c.Value = $n
}

function Set-NumCell($row, $col, $n) {
    if ($n -eq $null) { return }
    $c = $ws.Cells.Item($row, $col)
    $c.Font.Name = "Arial"
    $c.Font.Size = 9
    $c.NumberFormat = "0.00"
    $c.Value = $n
}

# ---- Header row (row 1) ----
Set-TextCell 1 1 "idx"
Set-TextCell 1 2 "idx2"
Set-TextCell 1 3 "Name"
Set-TextCell 1 4 "Date Start"
Set-TextCell 1 5 "Date End"
Set-TextCell 1 6 "(m3/s)"
Set-TextCell 1 7 "(MW1)"
Set-TextCell 1 8 "(MW2)"
Set-TextCell 1 9 "(GWh) Winter"
Set-TextCell 1 10 "(GWh) Summer"
Set-TextCell 1 11 "(GWh) Year"

# ---- Data rows 2-7, one per power plant, from the snapshot ----
$destRow = 2
foreach ($row in $data) {
    Set-IntCell $destRow 1 $row.A
    Set-IntCell $destRow 2 $row.B
    Set-TextCell $destRow 3 $row.C
    Set-IntCell $destRow 4 $row.D
    Set-IntCell $destRow 5 $row.E
    Set-NumCell $destRow 6 $row.F
    Set-NumCell $destRow 7 $row.G
    Set-NumCell $destRow 8 $row.H
    Set-NumCell $destRow 9 $row.I
    Set-NumCell $destRow 10 $row.J
    Set-NumCell $destRow 11 $row.K
    $destRow++
}

# ---- View state: selection moves to the new last data row ----
$ws.Range("A4:K4").Select()
